# Refresh the cryptos price list (rows 2-51) with the latest scraped
# values. Several "Price" (column D) cells look like plain numbers
# (e.g. "211.92", "1.00") but must stay as literal text so the exact
# formatting (trailing/leading zeros, etc.) is preserved instead of
# Excel silently reinterpreting them as numeric values. We force that
# by prefixing such values with a leading apostrophe (the same trick
# Excel's UI uses for "number stored as text").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.315.36'
$ws.Range("E2").Value = '  +0.38%  '
$ws.Range("D3").Value = '1.591.75'
$ws.Range("E3").Value = '  +0.62%  '
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").Value = '''211.92'
$ws.Range("E5").Value = '  +1.15%  '
$ws.Range("D6").Value = '''0.501'
$ws.Range("E6").Value = '  +0.23%  '
$ws.Range("E7").Value = '  -0.17%  '
$ws.Range("E8").Value = '  +0.27%  '
$ws.Range("E9").Value = '  -0.05%  '
$ws.Range("D10").Value = '''19.37'
$ws.Range("E10").Value = '  -0.75%  '
$ws.Range("D11").Value = '''0.0848'
$ws.Range("E11").Value = '  +0.33%  '
$ws.Range("D12").Value = '1.815.69'
$ws.Range("E12").Value = '  +0.71%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.579.90'
$ws.Range("E13").Value = '  +0.37%  '
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").Value = '''4.04'
$ws.Range("E14").Value = '  +0.38%  '
$ws.Range("D15").Value = '''0.521'
$ws.Range("E15").Value = '  +0.77%  '
$ws.Range("D16").Value = '''64.48'
$ws.Range("E16").Value = '  +0.12%  '
$ws.Range("D17").Value = '26.329.25'
$ws.Range("E17").Value = '  +0.47%  '
$ws.Range("D18").Value = '0.0₃0731'
$ws.Range("E18").Value = '  -0.38%  '
$ws.Range("D19").Value = '''7.49'
$ws.Range("E19").Value = '  +3.43%  '
$ws.Range("D20").Value = '''211.84'
$ws.Range("E20").Value = '  +2.74%  '
$ws.Range("E21").Value = '  -0.20%  '
$ws.Range("E22").Value = '  +0.85%  '
$ws.Range("E23").Value = '  +1.65%  '
$ws.Range("E24").Value = '  -2.24%  '
$ws.Range("D25").Value = '''145.30'
$ws.Range("E25").Value = '  +0.32%  '
$ws.Range("E26").Value = '  -0.05%  '
$ws.Range("E27").Value = '  +0.49%  '
$ws.Range("E28").Value = '  -0.58%  '
$ws.Range("D29").Value = '''15.20'
$ws.Range("E29").Value = '  -0.19%  '
$ws.Range("E30").Value = '  -0.26%  '
$ws.Range("E31").Value = '  +0.83%  '
$ws.Range("E32").Value = '  +0.11%  '
$ws.Range("E33").Value = '  +0.98%  '
$ws.Range("D34").Value = '1.335.24'
$ws.Range("E34").Value = '  +4.05%  '
$ws.Range("E35").Value = '  -0.96%  '
$ws.Range("E36").Value = '  +0.06%  '
$ws.Range("E37").Value = '  +0.35%  '
$ws.Range("E38").Value = '  +0.30%  '
$ws.Range("D39").Value = '''1.06'
$ws.Range("E39").Value = '  -15.19%  '
$ws.Range("D40").Value = '''0.819'
$ws.Range("D41").Value = '''5.78'
$ws.Range("E41").Value = '  +4.26%  '
$ws.Range("E42").Value = '  -0.15%  '
$ws.Range("E43").Value = '  +0.78%  '
$ws.Range("E44").Value = '  -0.82%  '
$ws.Range("D45").Value = '1.728.03'
$ws.Range("E45").Value = '  +0.63%  '
$ws.Range("E46").Value = '  -0.39%  '
$ws.Range("D47").Value = '''87.92'
$ws.Range("E47").Value = '  -0.76%  '
$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").Value = '''1.50'
$ws.Range("E48").Value = '  -3.33%  '
$ws.Range("B49").Value = 'Algorand'
$ws.Range("C49").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D49").Value = '''0.0982'
$ws.Range("E49").Value = '  -2.70%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").Value = '''0.0504'
$ws.Range("E50").Value = '  -0.55%  '
$ws.Range("B51").Value = 'USDD'
$ws.Range("C51").Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range("D51").Value = '''1.00'
$ws.Range("E51").Value = '  -0.36%  '
